$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4952
$ws.Range("G2").Value = 108
$ws.Range("F4").Value = 1867
$ws.Range("F5").Value = 2599
$ws.Range("F6").Value = 37
$ws.Range("F9").Value = 240
$ws.Range("F11").Value = 1094
$ws.Range("F12").Value = 363
$ws.Range("F14").Value = 59
$ws.Range("F16").Value = 16
$ws.Range("F17").Value = 234
$ws.Range("F18").Value = 128
$ws.Range("F19").Value = 80
$ws.Range("F20").Value = 1206
$ws.Range("F21").Value = 468
$ws.Range("F22").Value = 155
$ws.Range("F25").Value = 558
$ws.Range("F27").Value = 53
$ws.Range("F28").Value = 1924
$ws.Range("F29").Value = 2379
$ws.Range("F32").Value = 91
$ws.Range("F33").Value = 322
$ws.Range("F34").Value = 366
$ws.Range("F35").Value = 716
$ws.Range("F36").Value = 672
$ws.Range("F37").Value = 92
$ws.Range("F39").Value = 715
$ws.Range("F40").Value = 133
$ws.Range("F41").Value = 523
$ws.Range("F42").Value = 590
$ws.Range("F43").Value = 267
$ws.Range("F44").Value = 185

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 165
$ws.Range("F15").Value = 222
$ws.Range("F20").Value = 6

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 841

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 841
$ws.Range("F3").Value = 4952
$ws.Range("G3").Value = 108
$ws.Range("F4").Value = 1867
$ws.Range("F5").Value = 37
$ws.Range("F11").Value = 240
$ws.Range("F13").Value = 165
$ws.Range("F14").Value = 1094
$ws.Range("F15").Value = 363
$ws.Range("F17").Value = 59
$ws.Range("F19").Value = 234
$ws.Range("F21").Value = 128
$ws.Range("F22").Value = 80
$ws.Range("F23").Value = 1206
$ws.Range("F24").Value = 468
$ws.Range("F25").Value = 155
$ws.Range("F29").Value = 1924
$ws.Range("F30").Value = 2379
$ws.Range("F36").Value = 91
$ws.Range("F37").Value = 322
$ws.Range("F38").Value = 366
$ws.Range("F39").Value = 6
$ws.Range("F41").Value = 716
$ws.Range("F42").Value = 672
$ws.Range("F43").Value = 715
$ws.Range("F44").Value = 133
$ws.Range("F45").Value = 523
$ws.Range("F46").Value = 590
$ws.Range("F47").Value = 267
$ws.Range("F48").Value = 185
